$wb = $excel.ActiveWorkbook

# --- Fix the "cim" worksheet: correct the git target typo and remove the
# stray blank spacer rows that separated each process record. ---
$cim = $wb.Worksheets.Item("cim")
$cim.Range("C3").Value = "/web/<progDir>/cim"

# Delete the blank rows from the bottom up so row numbers above the
# deletion point stay valid while we work.
$cim.Rows("9:9").Delete()
$cim.Rows("7:7").Delete()
$cim.Rows("5:5").Delete()

# --- Fix the "pdfgen" worksheet: add the missing slash before "pdf". ---
$pdfgen = $wb.Worksheets.Item("pdfgen")
$pdfgen.Range("C3").Value = "/web/<progDir>/pdf"
